$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.178.18'
$ws.Range("E2").Value = '  +3.79%  '
$ws.Range("D3").Value = '1.785.96'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("D4").Value = '''0.9981'
$ws.Range("E4").Value = '  -0.51%  '
$ws.Range("D5").Value = '''336.16'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").Value = '''0.9954'
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("D7").Value = '''0.3830'
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("D8").Value = '''0.3439'
$ws.Range("E8").Value = '  +0.71%  '
$ws.Range("D9").Value = '''47.70'
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").Value = '''1.159'
$ws.Range("E10").Value = '  -2.55%  '
$ws.Range("D11").Value = '''0.07446'
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").Value = '''23.30'
$ws.Range("E12").Value = '  +7.26%  '
$ws.Range("D13").Value = '''0.9949'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").Value = '''6.418'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = '1.785.31'
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = '''7.140'
$ws.Range("E16").Value = '  +0.58%  '
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("D18").Value = '''0.06643'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("D20").Value = '''0.9957'
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").Value = '''17.53'
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("D22").Value = '''6.449'
$ws.Range("E22").Value = '  -1.21%  '
$ws.Range("D23").Value = '28.176.33'
$ws.Range("E23").Value = '  +3.80%  '
$ws.Range("D24").Value = '''12.10'
$ws.Range("E24").Value = '  -1.23%  '
$ws.Range("D25").Value = '''2.381'
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").Value = '''1.448'
$ws.Range("E26").Value = '  -0.55%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''20.93'
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("D28").Value = '''2.441'
$ws.Range("E28").Value = '  -2.34%  '
$ws.Range("D29").Value = '''153.86'
$ws.Range("E29").Value = '  -1.10%  '
$ws.Range("D30").Value = '1.986.78'
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("D31").Value = '''134.68'
$ws.Range("E31").Value = '  +0.62%  '
$ws.Range("D32").Value = '''6.175'
$ws.Range("E32").Value = '  +2.02%  '
$ws.Range("D33").Value = '''3.953'
$ws.Range("E33").Value = '  -1.07%  '
$ws.Range("D34").Value = '''0.08804'
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("E35").Value = '  -2.11%  '
$ws.Range("D36").Value = '''0.02432'
$ws.Range("E36").Value = '  +4.92%  '
$ws.Range("D37").Value = '''0.6870'
$ws.Range("E37").Value = '  +0.43%  '
$ws.Range("D38").Value = '''5.347'
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("D39").Value = '''0.06344'
$ws.Range("E39").Value = '  +0.76%  '
$ws.Range("D40").Value = '''0.2188'
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("D41").Value = '''1.513'
$ws.Range("E41").Value = '  -6.76%  '
$ws.Range("D42").Value = '''1.246'
$ws.Range("E42").Value = '  +0.94%  '
$ws.Range("D43").Value = '''8.379'
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("D44").Value = '''14.21'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '''0.9954'
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("D47").Value = '''3.850'
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("D48").Value = '''132.43'
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("E49").Value = '  -1.50%  '
$ws.Range("D50").Value = '''0.07435'
$ws.Range("E50").Value = '  +4.73%  '
$ws.Range("D51").Value = '''1.272'
$ws.Range("E51").Value = '  +7.76%  '
